$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New handed-back file for this run: d142da33-fb2d-4b61-88f4-a34802ac2cbf.md
# Adds one row to each of the three tables (Overview, zh-cn, de-de).
# ---------------------------------------------------------------------

$newMd   = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.md"
$newPath = "e2e\d142da33-fb2d-4b61-88f4-a34802ac2cbf.md"
$zhXlf   = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.c202557b0636e130a850b77a13065b8c95c91b66.zh-cn.xlf"
$deXlf   = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.c202557b0636e130a850b77a13065b8c95c91b66.de-de.xlf"

# =======================================================================
# Overview sheet
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $newMd
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md", "", "", $newPath) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-22 18:47:51"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# zh-cn sheet
# =======================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md", "", "", $newMd) | Out-Null
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = "2016-08-22 18:47:46"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md", "", "", $newMd) | Out-Null
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = "2016-08-22 18:48:13"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

# =======================================================================
# de-de sheet
# =======================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md", "", "", $newMd) | Out-Null
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = "2016-08-22 18:47:51"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md", "", "", $newMd) | Out-Null
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = "2016-08-22 18:48:21"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"
